$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Post-Press Cut / Final Trim" activity (previously on row 8) is moved up
# to row 5 (AShot verified for premier press); rows 5-7 shift down to 6-8.

$ws.Range("B5").Value = "Post-Press Cut"
$ws.Range("D5").Value = "Final Trim"
$ws.Range("G5").Value = "17"
$ws.Range("L5").Value = "Final Trim"
$ws.Range("M5").Value = "Cut for Press`nFinal Trim"

$ws.Range("B6").Value = "Print F 2x0"
$ws.Range("D6").Value = "8CS 40`" (#9)"
$ws.Range("G6").Value = "182"
$ws.Range("L6").Value = "8CS 40`" (#9)"
$ws.Range("M6").Value = "8CS 40`" (#9)`n8CP 40`" (#9)"

$ws.Range("B7").Value = "Laminate"
$ws.Range("D7").Value = "Farmout 1"
$ws.Range("G7").Value = "156"
$ws.Range("L7").Value = "Farmout 1"
$ws.Range("M7").Value = "Farmout 1`nFarmout 2`nFarmout 3"

$ws.Range("B8").Value = "Print F (Varnish 1x0)"
$ws.Range("D8").Value = "8CS 40`" (#9)"
$ws.Range("G8").Value = "20"
$ws.Range("L8").Value = "8CS 40`" (#9)"
$ws.Range("M8").Value = "8CS 40`" (#9)`n8CP 40`" (#9)"
